$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "RUG343.fasta"
$ws.Range("B2").Value = -44057.220981686129

$ws.Range("A3").Value = "RUG585.fasta"
$ws.Range("B3").Value = 229487.08200628409

$ws.Range("A4").Value = "RUG774.fasta"
$ws.Range("B4").Value = 183045.94028116838
